# Update Name of Algo
# Apply updated KNN imputed values to the result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value  = 6.851999999999999
$ws.Range("D6").Value  = -7.388
$ws.Range("D7").Value  = -7.258999999999999
$ws.Range("B8").Value  = 6.359
$ws.Range("D8").Value  = -7.419
$ws.Range("E11").Value = 12.636
$ws.Range("A12").Value = -21.734
$ws.Range("B12").Value = 6.381
$ws.Range("B14").Value = 7.543000000000001
$ws.Range("E14").Value = 12.498
$ws.Range("D19").Value = -7.981
$ws.Range("E19").Value = 12.704
$ws.Range("D21").Value = -7.25
$ws.Range("E21").Value = 13.364
$ws.Range("B22").Value = 7.037999999999999
$ws.Range("D24").Value = -7.431999999999999
